$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 375
$ws.Range("J2").Value = 350
$ws.Range("L2").Value = 350
$ws.Range("N2").Value = -576
$ws.Range("H28").Value = 1319.2
$ws.Range("I28").Value = 1277.2307
$ws.Range("J28").Value = 1397.1428
$ws.Range("K28").Value = 1277.2307
$ws.Range("L28").Value = 1397.1428
$ws.Range("M28").Value = -792.2307000000001
$ws.Range("N28").Value = -2367.1428
$ws.Range("H43").Value = 1782.6666
$ws.Range("I43").Value = 780.8
$ws.Range("J43").Value = 2498.2856
$ws.Range("K43").Value = 780.8
$ws.Range("L43").Value = 2498.2856
$ws.Range("M43").Value = -711.8
$ws.Range("N43").Value = -2636.2856
$ws.Range("H94").Value = 6471.143
$ws.Range("I94").Value = 4554.727
$ws.Range("K94").Value = 4554.727
$ws.Range("M94").Value = -4103.727
$ws.Range("H96").Value = 1178.7693
$ws.Range("I96").Value = 390.45456
$ws.Range("J96").Value = 5514.5
$ws.Range("K96").Value = 1171.36368
$ws.Range("L96").Value = 16543.5
$ws.Range("M96").Value = 201.6363200000001
$ws.Range("N96").Value = -19289.5
$ws.Range("H100").Value = 1764
$ws.Range("I100").Value = 1945
$ws.Range("K100").Value = 1945
$ws.Range("M100").Value = -1404
$ws.Range("H112").Value = 3146.1765
$ws.Range("J112").Value = 3092.3333
$ws.Range("L112").Value = 9276.999899999999
$ws.Range("N112").Value = -11492.9999
$ws.Range("H129").Value = 1778.7954
$ws.Range("I129").Value = 1359.8182
$ws.Range("K129").Value = 4079.4546
$ws.Range("M129").Value = 920.5454
$ws.Range("H132").Value = 411749.78
$ws.Range("I132").Value = 527118.2
$ws.Range("K132").Value = 1581354.6
$ws.Range("M132").Value = -1578824.6
$ws.Range("H138").Value = 1919.38
$ws.Range("J138").Value = 1951.0105
$ws.Range("L138").Value = 5853.0315
$ws.Range("N138").Value = -16133.0315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3799.5557
$ws.Range("I45").Value = 3399.4
$ws.Range("K45").Value = 3399.4
$ws.Range("M45").Value = -3022.4
$ws.Range("H46").Value = 107575.5
$ws.Range("J46").Value = 113434
$ws.Range("L46").Value = 113434
$ws.Range("N46").Value = -114072
$ws.Range("H102").Value = 1856.8125
$ws.Range("I102").Value = 1609.2307
$ws.Range("K102").Value = 1609.2307
$ws.Range("M102").Value = 12.76929999999993
$ws.Range("H122").Value = 3382.348
$ws.Range("I122").Value = 2588.2942
$ws.Range("J122").Value = 5632.1665
$ws.Range("K122").Value = 7764.882599999999
$ws.Range("L122").Value = 16896.4995
$ws.Range("M122").Value = -5314.882599999999
$ws.Range("N122").Value = -21796.4995
$ws.Range("H132").Value = 3934.6562
$ws.Range("I132").Value = 1880.5435
$ws.Range("K132").Value = 5641.6305
$ws.Range("M132").Value = -3111.6305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3454.6765
$ws.Range("J20").Value = 2984
$ws.Range("L20").Value = 2984
$ws.Range("N20").Value = -3478
$ws.Range("H99").Value = 2270.7932
$ws.Range("I99").Value = 2455.15
$ws.Range("J99").Value = 1861.1111
$ws.Range("K99").Value = 2455.15
$ws.Range("L99").Value = 1861.1111
$ws.Range("M99").Value = -957.1500000000001
$ws.Range("N99").Value = -4857.1111
$ws.Range("H105").Value = 4010
$ws.Range("I105").Value = 4273.5713
$ws.Range("J105").Value = 3395
$ws.Range("K105").Value = 4273.5713
$ws.Range("L105").Value = 3395
$ws.Range("M105").Value = -2526.5713
$ws.Range("N105").Value = -6889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19610482
$ws.Range("I31").Value = 23810972
$ws.Range("J31").Value = 8199.888999999999
$ws.Range("K31").Value = 23810972
$ws.Range("L31").Value = 8199.888999999999
$ws.Range("M31").Value = -23810677
$ws.Range("N31").Value = -8789.888999999999
$ws.Range("H34").Value = 19610482
$ws.Range("I34").Value = 23810972
$ws.Range("J34").Value = 8199.888999999999
$ws.Range("K34").Value = 23810972
$ws.Range("L34").Value = 8199.888999999999
$ws.Range("M34").Value = -23810770
$ws.Range("N34").Value = -8603.888999999999
$ws.Range("H105").Value = 692
$ws.Range("I105").Value = 204.66667
$ws.Range("K105").Value = 204.66667
$ws.Range("M105").Value = 1542.33333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14932980
$ws.Range("I4").Value = 627332.3
$ws.Range("J4").Value = 54670890
$ws.Range("K4").Value = 1881996.9
$ws.Range("L4").Value = 164012670
$ws.Range("M4").Value = -1881884.9
$ws.Range("N4").Value = -164012894
$ws.Range("H39").Value = 3794.8333
$ws.Range("J39").Value = 4081
$ws.Range("L39").Value = 12243
$ws.Range("N39").Value = -12831
$ws.Range("H55").Value = 5231
$ws.Range("J55").Value = 5231
$ws.Range("L55").Value = 15693
$ws.Range("N55").Value = -16047
$ws.Range("H86").Value = 449.5
$ws.Range("J86").Value = 499
$ws.Range("L86").Value = 1497
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 449.5
$ws.Range("J89").Value = 499
$ws.Range("L89").Value = 4491
$ws.Range("N89").Value = -16347
$ws.Range("H107").Value = 415.32257
$ws.Range("I107").Value = 217.76923
$ws.Range("J107").Value = 558
$ws.Range("K107").Value = 653.30769
$ws.Range("L107").Value = 1674
$ws.Range("M107").Value = 1266.69231
$ws.Range("N107").Value = -5514
$ws.Range("H141").Value = 6307.875
$ws.Range("I141").Value = 3653.8462
$ws.Range("J141").Value = 9444.454
$ws.Range("K141").Value = 10961.5386
$ws.Range("L141").Value = 28333.362
$ws.Range("M141").Value = -5781.5386
$ws.Range("N141").Value = -38693.362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 259.73334
$ws.Range("I2").Value = 119
$ws.Range("J2").Value = 420.57144
$ws.Range("K2").Value = 119
$ws.Range("L2").Value = 420.57144
$ws.Range("M2").Value = -6
$ws.Range("N2").Value = -646.5714399999999
$ws.Range("H70").Value = 6648.091
$ws.Range("I70").Value = 7383.3335
$ws.Range("J70").Value = 6372.375
$ws.Range("K70").Value = 7383.3335
$ws.Range("L70").Value = 6372.375
$ws.Range("M70").Value = -7113.3335
$ws.Range("N70").Value = -6912.375
$ws.Range("H73").Value = 6648.091
$ws.Range("I73").Value = 7383.3335
$ws.Range("J73").Value = 6372.375
$ws.Range("K73").Value = 7383.3335
$ws.Range("L73").Value = 7383.3335
$ws.Range("M73").Value = -6447.3335
$ws.Range("N73").Value = -8244.375
$ws.Range("H102").Value = 14171871
$ws.Range("J102").Value = 3766
$ws.Range("L102").Value = 3766
$ws.Range("N102").Value = -7010

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1156.8334
$ws.Range("I22").Value = 913.2857
$ws.Range("K22").Value = 913.2857
$ws.Range("M22").Value = -618.2857
$ws.Range("H27").Value = 1156.8334
$ws.Range("I27").Value = 913.2857
$ws.Range("K27").Value = 913.2857
$ws.Range("M27").Value = -806.2857
$ws.Range("H40").Value = 38469156
$ws.Range("I40").Value = 55559220
$ws.Range("J40").Value = 33342134
$ws.Range("K40").Value = 55559220
$ws.Range("L40").Value = 33342134
$ws.Range("M40").Value = -55559084
$ws.Range("N40").Value = -33342406
$ws.Range("H93").Value = 1187.8788
$ws.Range("I93").Value = 1078
$ws.Range("J93").Value = 1480.8889
$ws.Range("K93").Value = 1078
$ws.Range("L93").Value = 1480.8889
$ws.Range("M93").Value = 170
$ws.Range("N93").Value = -3976.8889
$ws.Range("H140").Value = 68773.53999999999
$ws.Range("J140").Value = 68773.53999999999
$ws.Range("L140").Value = 68773.53999999999
$ws.Range("N140").Value = -79133.53999999999
